$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Lejano Norte vs Centro Sur
$ws.Range("D2").Value = 2.166984222503732
$ws.Range("E2").Value = 9.527549690324223
$ws.Range("F2").Value = 0.09202912383054959
$ws.Range("G2").Value = 0.0001
$ws.Range("H2").Value = 0.0006000000000000001
$ws.Range("I2").Value = "**"

# Row 3: Lejano Norte vs Oeste Fria
$ws.Range("D3").Value = 3.358788406993761
$ws.Range("E3").Value = 15.036291887234292
$ws.Range("F3").Value = 0.2546957372586211
$ws.Range("G3").Value = 0.0001
$ws.Range("H3").Value = 0.0006000000000000001
$ws.Range("I3").Value = "**"

# Row 4: Lejano Norte vs Norte
$ws.Range("D4").Value = 1.896929866649466
$ws.Range("E4").Value = 10.102234669086721
$ws.Range("F4").Value = 0.23994533184492028
$ws.Range("G4").Value = 0.0001
$ws.Range("H4").Value = 0.0006000000000000001
$ws.Range("I4").Value = "**"

# Row 5: Centro Sur vs Oeste Fria
$ws.Range("D5").Value = 2.437830186568415
$ws.Range("E5").Value = 11.440913642772001
$ws.Range("F5").Value = 0.09268331953440556
$ws.Range("G5").Value = 0.0001
$ws.Range("H5").Value = 0.0006000000000000001
$ws.Range("I5").Value = "**"

# Row 6: Centro Sur vs Norte
$ws.Range("D6").Value = 0.7693646485375497
$ws.Range("E6").Value = 3.838179195521353
$ws.Range("F6").Value = 0.03696308260851031
$ws.Range("G6").Value = 0.002
$ws.Range("H6").Value = 0.012
$ws.Range("I6").Value = "."

# Row 7: Oeste Fria vs Norte
$ws.Range("D7").Value = 2.72331607196683
$ws.Range("E7").Value = 16.030878241766516
$ws.Range("F7").Value = 0.24277851012477528
$ws.Range("G7").Value = 0.0001
$ws.Range("H7").Value = 0.0006000000000000001
$ws.Range("I7").Value = "**"
